{"js": "const replacements = [\n  [\"2025-02-17 Monday\", \"2025-02-18 Tuesday\"],\n  [\"26\u00f73=\", \"18\u00f77=\"],\n  [\"12\u00f73=\", \"13\u00f74=\"],\n  [\"53\u00f74=\", \"43\u00f78=\"],\n  [\"13\u00f75=\", \"42\u00f73=\"],\n  [\"41\u00f77=\", \"53\u00f73=\"],\n  [\"98\u00f77=\", \"85\u00f74=\"],\n  [\"90\u00f77=\", \"32\u00f76=\"],\n  [\"55\u00f73=\", \"36\u00f72=\"],\n  [\"94\u00f76=\", \"43\u00f72=\"],\n  [\"49\u00f78=\", \"75\u00f75=\"],\n  [\"92\u00f74=\", \"27\u00f77=\"],\n  [\"73\u00f72=\", \"76\u00f76=\"],\n  [\"85\u00f77=\", \"22\u00f76=\"],\n  [\"34\u00f73=\", \"12\u00f75=\"],\n  [\"58\u00f79=\", \"23\u00f72=\"],\n  [\"99\u00f77=\", \"95\u00f76=\"],\n  [\"76\u00f78=\", \"40\u00f72=\"],\n  [\"13\u00f77=\", \"31\u00f78=\"],\n  [\"99\u00f73=\", \"70\u00f72=\"],\n  [\"40\u00f76=\", \"21\u00f72=\"],\n  [\"88\u00f78=\", \"80\u00f78=\"],\n  [\"22\u00f72=\", \"43\u00f79=\"],\n  [\"74\u00f77=\", \"59\u00f77=\"],\n  [\"45\u00f75=\", \"10\u00f74=\"],\n  [\"75\u00f73=\", \"89\u00f79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-17 Monday\", \"2025-02-18 Tuesday\"),\n    @(\"26\u00f73=\", \"18\u00f77=\"),\n    @(\"12\u00f73=\", \"13\u00f74=\"),\n    @(\"53\u00f74=\", \"43\u00f78=\"),\n    @(\"13\u00f75=\", \"42\u00f73=\"),\n    @(\"41\u00f77=\", \"53\u00f73=\"),\n    @(\"98\u00f77=\", \"85\u00f74=\"),\n    @(\"90\u00f77=\", \"32\u00f76=\"),\n    @(\"55\u00f73=\", \"36\u00f72=\"),\n    @(\"94\u00f76=\", \"43\u00f72=\"),\n    @(\"49\u00f78=\", \"75\u00f75=\"),\n    @(\"92\u00f74=\", \"27\u00f77=\"),\n    @(\"73\u00f72=\", \"76\u00f76=\"),\n    @(\"85\u00f77=\", \"22\u00f76=\"),\n    @(\"34\u00f73=\", \"12\u00f75=\"),\n    @(\"58\u00f79=\", \"23\u00f72=\"),\n    @(\"99\u00f77=\", \"95\u00f76=\"),\n    @(\"76\u00f78=\", \"40\u00f72=\"),\n    @(\"13\u00f77=\", \"31\u00f78=\"),\n    @(\"99\u00f73=\", \"70\u00f72=\"),\n    @(\"40\u00f76=\", \"21\u00f72=\"),\n    @(\"88\u00f78=\", \"80\u00f78=\"),\n    @(\"22\u00f72=\", \"43\u00f79=\"),\n    @(\"74\u00f77=\", \"59\u00f77=\"),\n    @(\"45\u00f75=\", \"10\u00f74=\"),\n    @(\"75\u00f73=\", \"89\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
